# overlap_data.xlsx: split single "2 Solution Design" sheet parse into
# three sheets — rename the original and add two new ones for
# "17 Data Loss Prevention (DLP)" and "18 Proxy Requirement Document".

$wb = $excel.ActiveWorkbook

# --- Sheet 1: rename existing sheet -----------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "2 Solution Design_vs_2 Solution"

# --- Sheet 2: "17 Data Loss Prevention (DLP)_v" ------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "17 Data Loss Prevention (DLP)_v"

$ws2.Range("A1").Value = "Number"
$ws2.Range("B1").Value = "Name"

$dlpRows = @(
    @(17.1,  "Does your solutions provide DLP for all traffic traversing the proxy?`nIf yes, please provide details about architecture, functionality, integration, etc."),
    @(17.3,  "Describe how your DLP solution classifies data. Does it support automatic classification based on predefined templates, and can custom classification rules be created?"),
    @(17.4,  "What methods does your DLP solution use to identify sensitive data (e.g., pattern matching, keyword analysis, machine learning)? Provide details on the accuracy and false positive rates."),
    @(17.5,  "Explain how your solution performs content inspection on data in motion (network traffic) and data at rest (stored data). What protocols and file types are supported?"),
    @(17.6,  "How does your DLP solution handle encrypted data? Can it inspect encrypted traffic and perform tokenization or redaction of sensitive information?"),
    @(17.7,  "What incident response capabilities does your DLP solution provide? Describe the alerting, reporting, and remediation features available for detected data breaches."),
    @(17.8,  "How does your DLP solution integrate with existing security infrastructure, such as SIEM, CASB, and IAM solutions? Provide examples of supported integrations and APIs."),
    @(17.9,  "Does your DLP solution incorporate UEBA to detect anomalous behavior? How does it differentiate between legitimate and malicious activities?"),
    @(17.1,  "Does your DLP solution incorporate UEBA (User and Entity Behavior Analytics) to detect anomalous behavior? How does it differentiate between legitimate and malicious activities?  If yes, indicate whether UEBA is built-in or using external UEBA information and the process to which that is configured."),
    @(17.12, "How does your solution monitor and control data movement across the network, including uploads to cloud storage, email attachments, and other outbound traffic?"),
    @(17.13, "How does your DLP solution assist in meeting regulatory compliance requirements (e.g., GDPR, HIPAA, PCI-DSS)? Provide details on compliance reporting and audit features."),
    @(17.14, "Describe the endpoint DLP capabilities of your solution. Can it monitor and control data transfer via USB drives, local storage, and other peripheral devices?"),
    @(17.15, "Explain the scalability of your DLP solution. How does it handle high traffic volumes and large-scale deployments across distributed environments?"),
    @(17.16, "What are the data retention and archiving policies for incidents and logs within your DLP solution? Can retention policies be customized based on organizational needs?"),
    @(17.17, "Does your DLP solution include features for user training and awareness? How does it educate employees on data protection practices and policy adherence?")
)

$r = 2
foreach ($row in $dlpRows) {
    $ws2.Cells.Item($r, 1).Value = $row[0]
    $ws2.Cells.Item($r, 2).Value = $row[1]
    $r = $r + 1
}

$ws2.Range("A1:B1").Font.Bold = $true
$ws2.Range("A1:B1").HorizontalAlignment = -4108
$ws2.Range("A1:B1").VerticalAlignment = -4160
$ws2.Range("A1:B1").Borders.LineStyle = 1

# --- Sheet 3: "18 Proxy Requirement Document_v" ------------------------
$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = "18 Proxy Requirement Document_v"

$ws3.Range("A1").Value = "Wells Fargo's Cloud Proxy and Zero  Trust Network Access RFP"
$ws3.Range("B1").Value = "Unnamed: 1"
$ws3.Range("C1").Value = "Unnamed: 2"
$ws3.Range("D1").Value = "Unnamed: 3"
$ws3.Range("E1").Value = "Wells Fargo's Cloud Proxy and Zero Trust Network Access RFP"

$ws3.Range("A2").Value = "Supplier Name"
$ws3.Range("E2").Value = "Supplier Name"

$ws3.Range("A1:E1").Font.Bold = $true
$ws3.Range("A1:E1").HorizontalAlignment = -4108
$ws3.Range("A1:E1").VerticalAlignment = -4160
$ws3.Range("A1:E1").Borders.LineStyle = 1

# Leave the first sheet ("2 Solution Design_vs_2 Solution") active/selected.
$ws1.Activate()
